$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.478.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "'1.886.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'246.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'42.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("D9").Value = "'56.55"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.74%  "
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").Value = "'0.0751"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("D12").Value = "'0.0984"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").Value = "'14.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.28%  "
$ws.Range("D14").Value = "'0.793"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.31%  "
$ws.Range("D15").Value = "'2.164.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "'5.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "'1.890.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "'35.500.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "'73.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").Value = "'0.0₃0829"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "'246.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'13.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("D23").Value = "'5.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.13%  "
$ws.Range("D24").Value = "'2.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'2.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").Value = "'165.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("D29").Value = "'18.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").Value = "'0.128"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "'4.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.78%  "
$ws.Range("E32").Value = "  +4.22%  "
$ws.Range("D33").Value = "'4.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'1.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +19.28%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -16.74%  "
$ws.Range("D37").Value = "'0.851"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "'0.0741"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.02%  "
$ws.Range("D39").Value = "'1.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.04%  "
$ws.Range("D40").Value = "'0.0229"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.87%  "
$ws.Range("D41").Value = "'98.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").Value = "'14.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +19.86%  "
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").Value = "'1.309.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").Value = "'2.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").Value = "'6.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "'42.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.04%  "
